$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-16 (5 sending clusters x 3 target clusters).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value, K Receptor-expressing cells,
# L Receptor detection rate, M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value, P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity
$data = @(
    @("ECs", "Gnai2", "Agtr2", "ECs", 3, 1, 144.783305, 434.349915, 0.2430046335191003, 0.251012682214973, 1, 0.3333333333333333, 0.01476566666666667, 0.044297, 0.003369932031170864, 0.003537680284170835, 2.137822020528333, 19.240398184755, 0.0008189090982189532, 0.0008880026169487492),
    @("ECs", "Gnai2", "Agtr2", "FAPs", 3, 1, 144.783305, 434.349915, 0.2430046335191003, 0.251012682214973, 3, 1, 3.743532666666667, 11.230598, 0.8543773151546027, 0.8969064524470826, 542.0010318554633, 4878.00928669917, 0.2076176463561771, 0.2251348943246583),
    @("ECs", "Gnai2", "Agtr2", "MuSCs", 3, 1, 144.783305, 434.349915, 0.2430046335191003, 0.251012682214973, 2, 1, 0.6232934999999999, 1.246587, 0.1422527528142264, 0.09955586726874661, 90.24249291501748, 541.4549574901049, 0.03456807806470426, 0.02498978527336592),
    @("FAPs", "Gnai2", "Agtr2", "ECs", 3, 1, 82.24887099999999, 246.746613, 0.1380466950572427, 0.1425959278859072, 1, 0.3333333333333333, 0.01476566666666667, 0.044297, 0.003369932031170864, 0.003537680284170835, 1.214459412895667, 10.930134716061, 0.0004652079794706786, 0.0005044588026850201),
    @("FAPs", "Gnai2", "Agtr2", "FAPs", 3, 1, 82.24887099999999, 246.746613, 0.1380466950572427, 0.1425959278859072, 3, 1, 3.743532666666667, 11.230598, 0.8543773151546027, 0.8969064524470826, 307.9013353849526, 2771.112018464574, 0.1179439646889732, 0.127895207813549),
    @("FAPs", "Gnai2", "Agtr2", "MuSCs", 3, 1, 82.24887099999999, 246.746613, 0.1380466950572427, 0.1425959278859072, 2, 1, 0.6232934999999999, 1.246587, 0.1422527528142264, 0.09955586726874661, 51.26518667663849, 307.591120059831, 0.01963752238879884, 0.01419626126967314),
    @("Inflammatory-Mac", "Gnai2", "Agtr2", "ECs", 3, 1, 163.8590903333333, 491.577271, 0.2750214756820535, 0.284084617144743, 1, 0.3333333333333333, 0.01476566666666667, 0.044297, 0.003369932031170864, 0.003537680284170835, 2.419488708165223, 21.775398373487, 0.0009268036801608308, 0.001005000549109177),
    @("Inflammatory-Mac", "Gnai2", "Agtr2", "FAPs", 3, 1, 163.8590903333333, 491.577271, 0.2750214756820535, 0.284084617144743, 3, 1, 3.743532666666667, 11.230598, 0.8543773151546027, 0.8969064524470826, 613.4118573931177, 5520.706716538059, 0.2349721100030897, 0.2547973261580791),
    @("Inflammatory-Mac", "Gnai2", "Agtr2", "MuSCs", 3, 1, 163.8590903333333, 491.577271, 0.2750214756820535, 0.284084617144743, 2, 1, 0.6232934999999999, 1.246587, 0.1422527528142264, 0.09955586726874661, 102.1323059206795, 612.7938355240769, 0.03912256199880294, 0.02828229043755473),
    @("MuSCs", "Gnai2", "Agtr2", "ECs", 2, 1, 57.0238095, 114.047619, 0.09570889357312636, 0.06590860906562239, 1, 0.3333333333333333, 0.01476566666666667, 0.044297, 0.003369932031170864, 0.003537680284170835, 0.8419945631405, 5.051967378843, 0.0003225324661200017, 0.0002331635868485755),
    @("MuSCs", "Gnai2", "Agtr2", "FAPs", 2, 1, 57.0238095, 114.047619, 0.09570889357312636, 0.06590860906562239, 3, 1, 3.743532666666667, 11.230598, 0.8543773151546027, 0.8969064524470826, 213.470493641027, 1280.822961846162, 0.08177150752742531, 0.05911385674276901),
    @("MuSCs", "Gnai2", "Agtr2", "MuSCs", 2, 1, 57.0238095, 114.047619, 0.09570889357312636, 0.06590860906562239, 2, 1, 0.6232934999999999, 1.246587, 0.1422527528142264, 0.09955586726874661, 35.54256980658825, 142.170279226353, 0.01361485357958105, 0.006561588736004813),
    @("Resolving-Mac", "Gnai2", "Agtr2", "ECs", 3, 1, 147.8896333333333, 443.6689, 0.2482183021684772, 0.2563981636887546, 1, 0.3333333333333333, 0.01476566666666667, 0.044297, 0.003369932031170864, 0.003537680284170835, 2.183689029255556, 19.6532012633, 0.0008364788072003996, 0.0009070547285793137),
    @("Resolving-Mac", "Gnai2", "Agtr2", "FAPs", 3, 1, 147.8896333333333, 443.6689, 0.2482183021684772, 0.2563981636887546, 3, 1, 3.743532666666667, 11.230598, 0.8543773151546027, 0.8969064524470826, 553.6296734446889, 4982.6670610022, 0.2120720865789374, 0.2299651674080272),
    @("Resolving-Mac", "Gnai2", "Agtr2", "MuSCs", 3, 1, 147.8896333333333, 443.6689, 0.2482183021684772, 0.2563981636887546, 2, 1, 0.6232934999999999, 1.246587, 0.1422527528142264, 0.09955586726874661, 92.17864717405, 553.0718830443, 0.03530973678233935, 0.02552594155214802)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}
